$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C92:C189").NumberFormat = "yyyy/mm/dd"
